# Apply "Penalty Reward System" forecast tweaks.
$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": update MyForecast (column D) values ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$newForecast = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
}

foreach ($row in $newForecast.Keys) {
    $wsForecast.Range("D$row").Value = $newForecast[$row]
}

# --- Sheet "Summary": update derived forecast metrics ---
# These cells were originally stored as text (even the numeric-looking
# ones), so force text formatting before writing to avoid Excel silently
# re-typing them as numbers / dates.
$wsSummary = $wb.Worksheets.Item("Summary")

$summaryCells = @{
    "B9"  = "16"            # Total Forecast (16 Weeks)
    "B10" = "8"             # Total Forecast (8 Weeks)
    "B11" = "4"             # Total Forecast (4 Weeks)
    "B12" = "2"             # Max Forecast
    "B14" = "1"             # Min Forecast
    "B15" = "2024-12-15"    # Min Forecast Week
}

foreach ($addr in $summaryCells.Keys) {
    $cell = $wsSummary.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryCells[$addr]
}
